$d = $word.ActiveDocument

function Find-ParaIndex($startsWith) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        if ($d.Paragraphs($i).Range.Text.StartsWith($startsWith)) {
            return $i
        }
    }
    return -1
}

function Insert-ParaBefore($startsWith, $text) {
    $idx = Find-ParaIndex($startsWith)
    $target = $d.Paragraphs($idx)
    $r = $target.Range
    $r.Collapse(1)
    $r.InsertParagraphBefore()
    $newPara = $d.Paragraphs($idx)
    $newPara.Range.Text = $text
}

# ------------------------------------------------------------------
# 1) "( Servlet & JSP , Oracle)" - merge the split "Servlet" run (and
#    its surrounding proofErr spell-check markers) back into a single
#    run so the text reads as one contiguous run.
# ------------------------------------------------------------------
$servletOld = "( Servlet & JSP , Oracle)"
$find1 = $d.Content.Find
$find1.Execute($servletOld, $true, $false, $false, $false, $false, $true, 1, $false, $servletOld, 2) | Out-Null

# ------------------------------------------------------------------
# 2) Insert two new paragraphs ("Header.jsp ..." / "Footer.jsp ...")
#    right before the "HOME PAGE (Hyper Links)" paragraph.
# ------------------------------------------------------------------
Insert-ParaBefore "HOME PAGE (Hyper Links)" "Header.jsp  :  Project Title  <h1>"
Insert-ParaBefore "HOME PAGE (Hyper Links)" "Footer.jsp   : Created By Yourname <h2>"

# ------------------------------------------------------------------
# 3) "ADD MARKS (Marks Entry of 5 subjects along with Student
#    Rollno(Primary Key), Student Name) " - merge the split "Rollno"
#    run (and its proofErr markers) into the surrounding run, while
#    keeping the "ADD MARKS" run and the trailing " " run untouched.
#    A plain text Find/Replace across the run-boundary would merge
#    every same-formatted run in the paragraph (including "ADD MARKS"
#    and the trailing space), so instead we splice in literal
#    WordOpenXML for just the paragraph in question.
# ------------------------------------------------------------------
$addMarksIdx = Find-ParaIndex("ADD MARKS")
$addMarksPara = $d.Paragraphs($addMarksIdx)

$xmlFrag = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>ADD MARKS</w:t></w:r><w:r><w:t xml:space="preserve"> (Marks Entry of 5 subjects along with Student Rollno(Primary Key), Student Name)</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$addMarksPara.Range.InsertXML($xmlFrag) | Out-Null
